$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
try {
    $ph = $s.Shapes.AddPlaceholder(2)
    Write-Output "AddPlaceholder ok: $ph Name=$($ph.Name)"
} catch {
    Write-Output "ERROR: $_"
}
